# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" column (D) for the
# 2d5489b1-3d79-483c-a379-dd25fef59254 row (row 6) on each
# language sheet, reflecting a freshly generated handoff report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D6").Value = "2016-03-09 09:48:47"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D6").Value = "2016-03-09 09:48:50"
